$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: F1 gains wrap-text styling (style s=1 -> s=3) ---
$ws.Range("F1").WrapText = $true

# --- Row 3: F3 gains wrap-text styling (no style -> s=2); text unchanged ---
$ws.Range("F3").WrapText = $true

# --- Row 4: C4 problem statement gets a bold "Leetcode" suffix appended ---
$ws.Range("C4").Value = "find target in rotated sorted array .  Use O(logn) Leetcode"
$ws.Range("C4").Characters(52, 8).Font.Bold = $true
$ws.Rows(4).RowHeight = 30

# --- Row 5: C5 problem statement gets a bold "Leetcode" suffix appended ---
$ws.Range("C5").Value = "find minimum in rotated sorted array in logn time complexity. Leetcode"
$ws.Range("C5").Characters(63, 8).Font.Bold = $true
$ws.Rows(5).RowHeight = 90

# --- Row 6: no content change, only auto row-height grows ---
$ws.Rows(6).RowHeight = 210

# --- Row 8 / 9: only auto row-height grows ---
$ws.Rows(8).RowHeight = 60
$ws.Rows(9).RowHeight = 60

# --- Row 10: C10 gets bold "Leetcode" suffix ---
$ws.Range("C10").Value = "Merged Sorted Arrays. Leetcode"
$ws.Range("C10").Characters(23, 8).Font.Bold = $true
$ws.Rows(10).RowHeight = 60

# --- Row 12: new data added (date, problem statement, solution, who) ---
$src = $ws.Range("B3")
$dst = $ws.Range("B12")
$src.Copy()
$dst.PasteSpecial(-4122)
$dst.Value = 45650

$ws.Range("C12").Value = "Trapping Rain Water. Leetcode"
$ws.Range("C12").Characters(22, 8).Font.Bold = $true
$ws.Range("C12").WrapText = $true

# --- Row 10: G10 text changes (allocated after "Trapping Rain Water" below) ---
$ws.Range("G10").Value = "solved and submitted "

$ws.Range("F12").Value = "for n^2 complexity :  For each element find left max and right max and in result add (min(left max, right max)-arr[i])" + [char]10 + "for n complexity : maintain leftmax and rightmax array ." + [char]10 + "For leftmax array : max(lm[i-1], arr[i])" + [char]10 + " rightmax array : max(rm[i+1], arr[i])" + [char]10 + " in result add (min(leftmax[i], rightmax[i])-arr[i])"
$ws.Range("F12").WrapText = $true

$ws.Range("G12").Value = "solved and submitted "
$ws.Range("G12").WrapText = $true

$ws.Rows(12).RowHeight = 120

# --- Final selection / active cell to match the edited region ---
$ws.Range("F12").Select()
